$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "2 hrs?"
$ws.Range("C3").Value = "Data Structures -- use set to parse through to find attributes"

$ws.Range("C4").Select()
